$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 28, shifting the existing rows 28-41 down to 29-42.
$ws.Rows.Item(28).Insert()

# Populate the new row 28 with the new weekly record.
$ws.Cells.Item(28, 1).Value = 11
$ws.Cells.Item(28, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(28, 3).Value = "Bíobío"
$ws.Cells.Item(28, 4).Value = 44875
$ws.Cells.Item(28, 4).NumberFormat = $ws.Cells.Item(29, 4).NumberFormat
$ws.Cells.Item(28, 5).Value = 8
$ws.Cells.Item(28, 6).Value = 100112026
$ws.Cells.Item(28, 7).Value = "Haba"
$ws.Cells.Item(28, 8).Value = "Sin especificar"
$ws.Cells.Item(28, 9).Value = "Primera"
$ws.Cells.Item(28, 10).Value = 200
$ws.Cells.Item(28, 11).Value = 6500
$ws.Cells.Item(28, 12).Value = 7000
$ws.Cells.Item(28, 13).Value = 6800
$ws.Cells.Item(28, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(28, 15).Value = "Región Metropolitana"
$ws.Cells.Item(28, 16).Value = 272
$ws.Cells.Item(28, 17).Value = 25
$ws.Cells.Item(28, 18).Value = "Hortaliza"

$wb.Save()
